$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update End Time (column E) from 100 to 80 for rows 2-55
$ws.Range("E2:E55").Value = 80

# Reflect the final selection left after the edit
$ws.Range("F2").Select()
